# Updated component input file given additional enclosure types.
# More comments for the script in case I die

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two existing insole enclosure entries to clarify they are
# "Manufactured" insoles (to distinguish from the new "Purchased" insole
# entries being added below).
$ws.Range("A33").Value = "2x Manufactured Insole with Removable Pod Enclosure"
$ws.Range("A34").Value = "2x Manufactured Insole with Attachable Pod Enclosure"

# Add the two new "Purchased" insole enclosure component rows.
$ws.Range("A35").Value = "2x Purchased Insole with Removable Pod Enclosure"
$ws.Range("B35").Value = "ENC_004"
$ws.Range("C35").Value = 3.9658
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 1

$ws.Range("A36").Value = "2x Purchased Insole with Attachable Pod Enclosure"
$ws.Range("B36").Value = "ENC_005"
$ws.Range("C36").Value = 4.3434
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 1


# Move the selection to reflect where editing left off.
$ws.Range("C37").Select()
